$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Set the header-row label cells (columns 2-12) to bold first.
for ($col = 2; $col -le 12; $col++) {
    $t.Cell(1, $col).Range.Font.Bold = $true
}

# Replace the label text in each header cell.
$d.Content.Find.Execute("Chiffre d’affaires (K€)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Revenu (K`$)", 2) | Out-Null

$d.Content.Find.Execute("Coût des marchandises vendues (`$K)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Coût des marchandises vendues (K`$)", 2) | Out-Null

$d.Content.Find.Execute("Marge bénéficiaire brute (%)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Marge bénéficiaire brut (%)", 2) | Out-Null

$d.Content.Find.Execute("Dépenses de fonctionnement (`$K)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Dépenses opérationnelles (K`$)", 2) | Out-Null

$d.Content.Find.Execute("EBITDA (`$K)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "EBITDA (K`$)", 2) | Out-Null

$d.Content.Find.Execute("Charges d’intérêt (`$K)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Charges d’intérêts (K`$)", 2) | Out-Null

$d.Content.Find.Execute("Bénéfice avant impôts (`$K)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Bénéfice avant impôt (K`$)", 2) | Out-Null

$d.Content.Find.Execute("Revenus nets (`$K)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Résultat net (K`$)", 2) | Out-Null

$d.Content.Find.Execute("Total des actifs (`$K)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Total actif (K`$)", 2) | Out-Null

$d.Content.Find.Execute("Total du passif (`$K)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Total passif (K`$)", 2) | Out-Null

$d.Content.Find.Execute("Capitaux propres (`$K)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Capitaux propres des actionnaires (K`$)", 2) | Out-Null
